# Auto-generated Excel COM-interop edit script
# Applies cached-value updates to the Leve profit calculation columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 67333.336
$ws.Range("I52").Value = 1000
$ws.Range("J52").Value = 200000
$ws.Range("K52").Value = 3000
$ws.Range("L52").Value = 600000
$ws.Range("M52").Value = -2840
$ws.Range("N52").Value = -600320
$ws.Range("H98").Value = 2445.16
$ws.Range("I98").Value = 1187.6818
$ws.Range("J98").Value = 11666.667
$ws.Range("K98").Value = 1187.6818
$ws.Range("L98").Value = 11666.667
$ws.Range("M98").Value = 310.3181999999999
$ws.Range("N98").Value = -14662.667
$ws.Range("H106").Value = 4214
$ws.Range("I106").Value = 2916.3333
$ws.Range("K106").Value = 2916.3333
$ws.Range("M106").Value = -2285.3333
$ws.Range("H111").Value = 2000
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 2000
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 6000
$ws.Range("M111").ClearContents()
$ws.Range("N111").Value = -12134
$ws.Range("H122").Value = 2445.16
$ws.Range("I122").Value = 1187.6818
$ws.Range("J122").Value = 11666.667
$ws.Range("K122").Value = 3563.0454
$ws.Range("L122").Value = 35000.001
$ws.Range("M122").Value = -1113.0454
$ws.Range("N122").Value = -39900.001
$ws.Range("H125").Value = 3488.889
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 3488.889
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 31400.001
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -36320.001
$ws.Range("H129").Value = 918.7659
$ws.Range("J129").Value = 963.5814
$ws.Range("L129").Value = 2890.7442
$ws.Range("N129").Value = -12890.7442
$ws.Range("H132").Value = 21743628
$ws.Range("I132").Value = 32263222
$ws.Range("J132").Value = 3131.8
$ws.Range("K132").Value = 96789666
$ws.Range("L132").Value = 9395.400000000001
$ws.Range("M132").Value = -96787136
$ws.Range("N132").Value = -14455.4
$ws.Range("H135").Value = 839.2857
$ws.Range("I135").Value = 395.1
$ws.Range("J135").Value = 1949.75
$ws.Range("K135").Value = 3555.9
$ws.Range("L135").Value = 17547.75
$ws.Range("M135").Value = -1020.9
$ws.Range("N135").Value = -22617.75
$ws.Range("H137").Value = 1782.6613
$ws.Range("I137").Value = 897.1836499999999
$ws.Range("J137").Value = 5120.231
$ws.Range("K137").Value = 2691.55095
$ws.Range("L137").Value = 15360.693
$ws.Range("M137").Value = -141.5509499999998
$ws.Range("N137").Value = -20460.693
$ws.Range("H138").Value = 2727.6
$ws.Range("I138").Value = 811.2432
$ws.Range("J138").Value = 3853.0793
$ws.Range("K138").Value = 2433.7296
$ws.Range("L138").Value = 11559.2379
$ws.Range("M138").Value = 2706.2704
$ws.Range("N138").Value = -21839.2379
$ws.Range("H141").Value = 7833.9697
$ws.Range("I141").Value = 8246.241
$ws.Range("J141").Value = 4845
$ws.Range("K141").Value = 24738.723
$ws.Range("L141").Value = 14535
$ws.Range("M141").Value = -19558.723
$ws.Range("N141").Value = -24895

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3583.484
$ws.Range("I32").Value = 2545.2942
$ws.Range("J32").Value = 14614.25
$ws.Range("K32").Value = 2545.2942
$ws.Range("L32").Value = 14614.25
$ws.Range("M32").Value = -2258.2942
$ws.Range("N32").Value = -15188.25
$ws.Range("H61").Value = 852.6
$ws.Range("I61").Value = 723.5333000000001
$ws.Range("J61").Value = 1110.7333
$ws.Range("K61").Value = 723.5333000000001
$ws.Range("L61").Value = 1110.7333
$ws.Range("M61").Value = -511.5333000000001
$ws.Range("N61").Value = -1534.7333
$ws.Range("H63").Value = 9895658
$ws.Range("I63").Value = 10656786
$ws.Range("J63").Value = 1004
$ws.Range("K63").Value = 10656786
$ws.Range("L63").Value = 1004
$ws.Range("M63").Value = -10656100
$ws.Range("N63").Value = -2376
$ws.Range("H66").Value = 9895658
$ws.Range("I66").Value = 10656786
$ws.Range("J66").Value = 1004
$ws.Range("K66").Value = 53283930
$ws.Range("L66").Value = 5020
$ws.Range("M66").Value = -53280498
$ws.Range("N66").Value = -11884
$ws.Range("H74").Value = 2300.2593
$ws.Range("I74").Value = 2219.6086
$ws.Range("J74").Value = 2764
$ws.Range("K74").Value = 2219.6086
$ws.Range("L74").Value = 2764
$ws.Range("M74").Value = -1345.6086
$ws.Range("N74").Value = -4512
$ws.Range("H77").Value = 2300.2593
$ws.Range("I77").Value = 2219.6086
$ws.Range("J77").Value = 2764
$ws.Range("K77").Value = 11098.043
$ws.Range("L77").Value = 13820
$ws.Range("M77").Value = -6730.043
$ws.Range("N77").Value = -22556
$ws.Range("H132").Value = 1795.6792
$ws.Range("I132").Value = 1251.475
$ws.Range("J132").Value = 3470.1538
$ws.Range("K132").Value = 3754.425
$ws.Range("L132").Value = 10410.4614
$ws.Range("M132").Value = -1224.425
$ws.Range("N132").Value = -15470.4614
$ws.Range("H136").Value = 852.6
$ws.Range("I136").Value = 723.5333000000001
$ws.Range("J136").Value = 1110.7333
$ws.Range("K136").Value = 2170.5999
$ws.Range("L136").Value = 3332.199900000001
$ws.Range("M136").Value = 379.4000999999998
$ws.Range("N136").Value = -8432.1999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3048.0312
$ws.Range("I134").Value = 1497.2273
$ws.Range("J134").Value = 6459.8
$ws.Range("K134").Value = 4491.6819
$ws.Range("L134").Value = 19379.4
$ws.Range("M134").Value = -1956.6819
$ws.Range("N134").Value = -24449.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9261503
$ws.Range("I31").Value = 1179.1892
$ws.Range("K31").Value = 1179.1892
$ws.Range("M31").Value = -884.1892
$ws.Range("H34").Value = 9261503
$ws.Range("I34").Value = 1179.1892
$ws.Range("K34").Value = 1179.1892
$ws.Range("M34").Value = -977.1892
$ws.Range("H58").Value = 1632.2471
$ws.Range("I58").Value = 1413.1688
$ws.Range("J58").Value = 3740.875
$ws.Range("K58").Value = 1413.1688
$ws.Range("L58").Value = 3740.875
$ws.Range("M58").Value = -1210.1688
$ws.Range("N58").Value = -4146.875
$ws.Range("H122").Value = 2387.077
$ws.Range("I122").Value = 1133.1428
$ws.Range("K122").Value = 3399.4284
$ws.Range("M122").Value = -949.4284000000002
$ws.Range("H132").Value = 2348.6538
$ws.Range("I132").Value = 1971.326
$ws.Range("J132").Value = 5241.5
$ws.Range("K132").Value = 5913.978
$ws.Range("L132").Value = 15724.5
$ws.Range("M132").Value = -3383.978
$ws.Range("N132").Value = -20784.5
$ws.Range("H134").Value = 3982.5227
$ws.Range("I134").Value = 4505.2144
$ws.Range("K134").Value = 13515.6432
$ws.Range("M134").Value = -10980.6432
$ws.Range("H136").Value = 1632.2471
$ws.Range("I136").Value = 1413.1688
$ws.Range("J136").Value = 3740.875
$ws.Range("K136").Value = 4239.5064
$ws.Range("L136").Value = 11222.625
$ws.Range("M136").Value = -1689.5064
$ws.Range("N136").Value = -16322.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 437.8
$ws.Range("I118").Value = 437.8
$ws.Range("K118").Value = 1313.4
$ws.Range("M118").Value = -70.40000000000009
$ws.Range("H140").Value = 2393.7097
$ws.Range("I140").Value = 2533.611
$ws.Range("J140").Value = 2200
$ws.Range("K140").Value = 7600.833
$ws.Range("L140").Value = 6600
$ws.Range("M140").Value = -2420.833
$ws.Range("N140").Value = -16960

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 40000
$ws.Range("J95").Value = 40000
$ws.Range("L95").Value = 40000
$ws.Range("N95").Value = -45492
$ws.Range("H126").Value = 2666.13
$ws.Range("I126").Value = 2666.13
$ws.Range("K126").Value = 7998.39
$ws.Range("M126").Value = -5528.39
$ws.Range("H132").Value = 2387.634
$ws.Range("I132").Value = 1383.8
$ws.Range("J132").Value = 5125.364
$ws.Range("K132").Value = 4151.4
$ws.Range("L132").Value = 15376.092
$ws.Range("M132").Value = -1621.4
$ws.Range("N132").Value = -20436.092

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 29418638
$ws.Range("I40").Value = 71433560
$ws.Range("J40").Value = 8194.9
$ws.Range("K40").Value = 71433560
$ws.Range("L40").Value = 8194.9
$ws.Range("M40").Value = -71433424
$ws.Range("N40").Value = -8466.9
$ws.Range("H132").Value = 7272.3687
$ws.Range("I132").Value = 7046.9756
$ws.Range("J132").Value = 7849.9375
$ws.Range("K132").Value = 21140.9268
$ws.Range("L132").Value = 23549.8125
$ws.Range("M132").Value = -18610.9268
$ws.Range("N132").Value = -28609.8125
$ws.Range("H136").Value = 1977.0702
$ws.Range("I136").Value = 992.907
$ws.Range("J136").Value = 4999.857
$ws.Range("K136").Value = 2978.721
$ws.Range("L136").Value = 14999.571
$ws.Range("M136").Value = -428.721
$ws.Range("N136").Value = -20099.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 22728864
$ws.Range("I81").Value = 22728864
$ws.Range("K81").Value = 45457728
$ws.Range("M81").Value = -45456667
$ws.Range("H84").Value = 22728864
$ws.Range("I84").Value = 22728864
$ws.Range("K84").Value = 227288640
$ws.Range("M84").Value = -227283336
$ws.Range("H92").Value = 29666.666
$ws.Range("J92").Value = 29666.666
$ws.Range("L92").Value = 29666.666
$ws.Range("N92").Value = -34658.666
$ws.Range("H132").Value = 5377733
$ws.Range("I132").Value = 964.06384
$ws.Range("K132").Value = 2892.19152
$ws.Range("M132").Value = -362.1915200000003
$ws.Range("H136").Value = 2183.182
$ws.Range("I136").Value = 602
$ws.Range("J136").Value = 5571.4287
$ws.Range("K136").Value = 1806
$ws.Range("L136").Value = 16714.2861
$ws.Range("M136").Value = 744
$ws.Range("N136").Value = -21814.2861
